$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 00:16"

# Re-sort: country names shift to new rows (same totals-B ranking, new day's data)
# Column A (country names) updates - repoints to existing shared strings
$ws.Cells.Item(115, 1).Value = "Ghana"
$ws.Cells.Item(116, 1).Value = "Bolivia"
$ws.Cells.Item(121, 1).Value = "Puerto Rico"
$ws.Cells.Item(122, 1).Value = "Montenegro"
$ws.Cells.Item(123, 1).Value = "Guam"
$ws.Cells.Item(127, 1).Value = "Guatemala"
$ws.Cells.Item(128, 1).Value = "Nigeria"
$ws.Cells.Item(129, 1).Value = "Monaco"
$ws.Cells.Item(130, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(132, 1).Value = "Etiopia"
$ws.Cells.Item(133, 1).Value = "Togo"
$ws.Cells.Item(134, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(136, 1).Value = "Seychelles"
$ws.Cells.Item(137, 1).Value = "Kenia"
$ws.Cells.Item(139, 1).Value = "Kirguistan"
$ws.Cells.Item(140, 1).Value = "Tanzania"
$ws.Cells.Item(141, 1).Value = "Mayotte"
$ws.Cells.Item(142, 1).Value = "Mongolia"
$ws.Cells.Item(145, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(146, 1).Value = "Surinam"
$ws.Cells.Item(147, 1).Value = "Bahamas"
$ws.Cells.Item(149, 1).Value = "San Bartolome"
$ws.Cells.Item(150, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(151, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(152, 1).Value = "Congo"
$ws.Cells.Item(153, 1).Value = "Madagascar"
$ws.Cells.Item(155, 1).Value = "Islas Caimanes"
$ws.Cells.Item(156, 1).Value = "Curazao"
$ws.Cells.Item(157, 1).Value = "Guinea"
$ws.Cells.Item(158, 1).Value = "Santa Lucia"
$ws.Cells.Item(161, 1).Value = "Butan"
$ws.Cells.Item(162, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(163, 1).Value = "Bermudas"
$ws.Cells.Item(164, 1).Value = "Liberia"
$ws.Cells.Item(165, 1).Value = "Nicaragua"
$ws.Cells.Item(166, 1).Value = "Benin"
$ws.Cells.Item(172, 1).Value = "Suazilandia"
$ws.Cells.Item(173, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(175, 1).Value = "El Salvador"
$ws.Cells.Item(176, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(177, 1).Value = "Zimbabue"
$ws.Cells.Item(178, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(179, 1).Value = "Montserrat"
$ws.Cells.Item(180, 1).Value = "Cabo Verde"
$ws.Cells.Item(181, 1).Value = "Gambia"
$ws.Cells.Item(182, 1).Value = "Republica del Chad"
$ws.Cells.Item(183, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(184, 1).Value = "Angola"
$ws.Cells.Item(185, 1).Value = "Fiyi"
$ws.Cells.Item(186, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(187, 1).Value = "Niger"

# Column B/C/D/E/F/G/H numeric updates
$ws.Cells.Item(9, 2).Value = 19195
$ws.Cells.Item(9, 3).Value = 5406
$ws.Cells.Item(9, 5).Value = 18799
$ws.Cells.Item(9, 7).Value = 42
$ws.Cells.Item(9, 8).Value = 249
$ws.Cells.Item(115, 3).Value = 5
$ws.Cells.Item(116, 3).Value = 1
$ws.Cells.Item(121, 3).Value = 8
$ws.Cells.Item(122, 3).Value = 1
$ws.Cells.Item(123, 3).Value = 2
$ws.Cells.Item(127, 3).Value = 3
$ws.Cells.Item(127, 4).Value = 0
$ws.Cells.Item(127, 8).Value = 1
$ws.Cells.Item(128, 3).Value = 0
$ws.Cells.Item(128, 4).Value = 1
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(129, 3).Value = 1
$ws.Cells.Item(130, 3).Value = 5
$ws.Cells.Item(132, 3).Value = 2
$ws.Cells.Item(133, 3).Value = 8
$ws.Cells.Item(134, 3).Value = 0
$ws.Cells.Item(136, 3).Value = 1
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(139, 3).Value = 3
$ws.Cells.Item(140, 3).Value = 0
$ws.Cells.Item(141, 3).Value = 2
$ws.Cells.Item(142, 3).Value = 0
$ws.Cells.Item(145, 3).Value = 1
$ws.Cells.Item(146, 3).Value = 3
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(150, 3).Value = 2
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(153, 3).Value = 3
$ws.Cells.Item(157, 3).Value = 1
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(165, 3).Value = 1
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(173, 3).Value = 1
$ws.Cells.Item(177, 3).Value = 1
$ws.Cells.Item(180, 3).Value = 1
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(187, 3).Value = 0
